# Upload excel files with prices
#
# 1) Refresh the scrape timestamp (column O) on every data row from
#    "2022-07-28 07:00:55" -> "2022-07-28 20:59:35".
# 2) Re-order the three "sock" rows (14-16): the row that used to be
#    last (id 6365813008, "Schwarz") now sorts first, and the other two
#    rows shift down by one (a cyclic rotation of rows 14,15,16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldStamp = "2022-07-28 07:00:55"
$newStamp = "2022-07-28 20:59:35"

# Helper: does a string look like something Excel would auto-parse as a
# number (so we must quote-prefix it to keep it text, matching the
# source data which stores ids/prices/etc. as text)?
function Test-LooksNumeric {
    param([string]$s)
    if ([string]::IsNullOrEmpty($s)) {
        return $false
    }
    if ($s -match '^-?\d+(\.\d+)?$') {
        return $true
    }
    return $false
}

# Helper: write a plain (non-numeric-column) cell back as text, even if
# its contents look numeric.
function Set-TextCell {
    param($cell, [string]$s)
    if ([string]::IsNullOrEmpty($s)) {
        $cell.Value = "'"
    } elseif (Test-LooksNumeric $s) {
        $cell.Value = "'" + $s
    } else {
        $cell.Value = $s
    }
}

$lastRow = 130
$lastCol = 15

# --- Step 1: bump the timestamp on every data row ------------------------
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $lastCol)
    if ($cell.Text -eq $oldStamp) {
        $cell.Value = $newStamp
    }
}

# --- Step 2: rotate rows 14, 15, 16 ---------------------------------------
# Capture the current contents (as text) of each row first, so the writes
# below can't clobber data we still need to read.
$capturedRows = @()
for ($r = 14; $r -le 16; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals += $ws.Cells.Item($r, $c).Text
    }
    $capturedRows += ,$rowVals
}

# capturedRows[0] = old row14, [1] = old row15, [2] = old row16
# new row14 = old row16, new row15 = old row14, new row16 = old row15
$newOrder = @($capturedRows[2], $capturedRows[0], $capturedRows[1])

# Columns E (5) and F (6) are the only numeric-typed columns in this sheet.
$numericCols = @(5, 6)

for ($i = 0; $i -lt 3; $i++) {
    $r = 14 + $i
    $rowVals = $newOrder[$i]
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $rowVals[$c - 1]
        if ($c -eq $lastCol) {
            # timestamp column -> always the fresh stamp
            $cell.Value = $newStamp
        } elseif ($numericCols -contains $c) {
            if ([string]::IsNullOrEmpty($val)) {
                $cell.Value = "'"
            } else {
                $cell.Value = [double]$val
            }
        } else {
            Set-TextCell $cell $val
        }
    }
}
